# Weekly update: insert two new "semanal" observation rows (Primera/Segunda)
# for Coliflor at the top of the historical data block, pushing the rest of
# the rows (old 1064:1174) down by two rows to (1066:1176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current first data row of this
# block (row 1064); this shifts all existing rows 1064-1174 down to 1066-1176,
# carrying their values/formatting with them (dimension grows to A1:R1176).
$ws.Rows("1064:1065").Insert()

# New row 1064: Primera, date 45194 (2023-09-25)
$ws.Cells.Item(1064, 1).Value = 3
$ws.Cells.Item(1064, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1064, 3).Value = "Coquimbo"
$ws.Cells.Item(1064, 4).Value = 45194
$ws.Cells.Item(1064, 5).Value = 5
$ws.Cells.Item(1064, 6).Value = 100112008
$ws.Cells.Item(1064, 7).Value = "Coliflor"
$ws.Cells.Item(1064, 8).Value = "Sin especificar"
$ws.Cells.Item(1064, 9).Value = "Primera"
$ws.Cells.Item(1064, 10).Value = 2300
$ws.Cells.Item(1064, 11).Value = 750
$ws.Cells.Item(1064, 12).Value = 800
$ws.Cells.Item(1064, 13).Value = 776
$ws.Cells.Item(1064, 14).Value = "`$/unidad"
$ws.Cells.Item(1064, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1064, 16).Value = 776
$ws.Cells.Item(1064, 17).Value = 1
$ws.Cells.Item(1064, 18).Value = "Hortaliza"

# New row 1065: Segunda, same date 45194
$ws.Cells.Item(1065, 1).Value = 3
$ws.Cells.Item(1065, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1065, 3).Value = "Coquimbo"
$ws.Cells.Item(1065, 4).Value = 45194
$ws.Cells.Item(1065, 5).Value = 5
$ws.Cells.Item(1065, 6).Value = 100112008
$ws.Cells.Item(1065, 7).Value = "Coliflor"
$ws.Cells.Item(1065, 8).Value = "Sin especificar"
$ws.Cells.Item(1065, 9).Value = "Segunda"
$ws.Cells.Item(1065, 10).Value = 1000
$ws.Cells.Item(1065, 11).Value = 600
$ws.Cells.Item(1065, 12).Value = 600
$ws.Cells.Item(1065, 13).Value = 600
$ws.Cells.Item(1065, 14).Value = "`$/unidad"
$ws.Cells.Item(1065, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1065, 16).Value = 600
$ws.Cells.Item(1065, 17).Value = 1
$ws.Cells.Item(1065, 18).Value = "Hortaliza"
